$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# New issue #8 (row 9) is now Fixed -> set its Status cell (D9), matching the
# same value/format used by the other Status cells in column D (e.g. D3, D7).
$ws.Range("D9").Value = "Fixed"
$ws.Range("D8").Copy()
$ws.Range("D9").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Scroll the sheet so column C / row 7 is the top-left visible cell, and move
# the active selection to the newly-filled D9 cell.
$win = $excel.ActiveWindow
$win.ScrollColumn = 3
$win.ScrollRow = 7
$ws.Range("D9").Select()
